$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the header row: the old "Testcases / Data1 / Data2" header becomes
# a "firstname / lastname / postcode" header.
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# The fourth ("...3") data column is dropped entirely - delete D1:D4 and
# shift the remaining cells left so the used range narrows back to A:C.
$ws.Range("D1:D4").Delete()

# Leave column D selected (whole-column selection), matching the resulting
# view state after clearing that column out.
$ws.Columns("D").Select() | Out-Null

# Widen/resize the workbook window.
$win = $excel.ActiveWindow
$win.Width = 10500
$win.Height = 6660
